$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '30.025.45'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  -0.89%  '

$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.902.27'
$ws.Range('D3').Style = 'Normal'

$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.000'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  -0.04%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '0.7415'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -1.21%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '243.60'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +0.50%  '

$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '1.000'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -0.01%  '

$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3067'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -3.56%  '

$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '26.13'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -5.79%  '

$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.06901'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -2.96%  '

$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.08028'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -0.11%  '

$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.7632'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -2.16%  '

$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '1.901.50'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -2.00%  '

$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '5.232'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -2.95%  '

$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '91.33'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -1.90%  '

$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '30.031.78'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -0.87%  '

$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '6.092'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +0.96%  '

$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '14.02'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -3.44%  '

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.000007754'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -2.33%  '

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '237.71'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -5.60%  '

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '1.001'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +0.10%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '2.157.45'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -1.00%  '

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '0.9999'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -0.07%  '

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '7.058'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +5.43%  '

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '9.319'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -2.29%  '

$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '166.38'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +0.94%  '

$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '18.82'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -1.48%  '

$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '0.1263'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -2.72%  '

$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '2.037'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -7.07%  '

$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.351'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -1.16%  '

$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.537'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -0.63%  '

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '4.289'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -2.79%  '

$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '4.038'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -2.47%  '

$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.05303'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +1.62%  '

$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.292'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -1.92%  '

$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.7363'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -2.63%  '

$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.708'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -2.78%  '

$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.01946'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -0.28%  '

$ws.Range('E39').Value = '  -0.19%  '

$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '6.271'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -3.50%  '

$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.4453'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -1.61%  '

$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '73.13'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -6.77%  '

$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '1.958'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -1.01%  '

$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '1.001'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +0.04%  '

$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.8329'
$ws.Range('D45').Style = 'Normal'

$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '7.624'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -0.74%  '

$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '101.33'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -0.35%  '

$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '9.814'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -2.27%  '

$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '2.057.49'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -0.99%  '

$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '36.49'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -3.51%  '

$ws.Range('B51').Value = 'Cronos'
$ws.Range('C51').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.05956'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -1.31%  '
